$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet's column definitions already carry the date / integer / 2-decimal
# styles (column A -> style 1, C:E -> style 3, F -> style 2), so simply
# writing values into row 36/37 picks up the same formatting as every other
# data row without needing an explicit copy/paste of row 34:35.

# Row 36: 四方坪站 (site string index 4)
$ws.Cells.Item(36, 1).Value = 45979
$ws.Cells.Item(36, 2).Value = "四方坪站"
$ws.Cells.Item(36, 3).Value = 9993.93
$ws.Cells.Item(36, 4).Value = 8862.26
$ws.Cells.Item(36, 5).Value = 3294.21
$ws.Cells.Item(36, 6).Value = 410

# Row 37: 高岭站 (site string index 5)
$ws.Cells.Item(37, 1).Value = 45979
$ws.Cells.Item(37, 2).Value = "高岭站"
$ws.Cells.Item(37, 3).Value = 5793.18
$ws.Cells.Item(37, 4).Value = 5175.38
$ws.Cells.Item(37, 5).Value = 1547.65
$ws.Cells.Item(37, 6).Value = 212

# Scroll the view down to match where the new rows were added and move the
# active selection to mirror the author's last edit position.
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("H36").Select()
